$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are plain text (e.g. "215.56", "27.172.93") in the source data.
# Excel auto-converts numeric-looking strings assigned via .Value to real numbers,
# so force text format, assign, then restore the default "Normal" style/format
# (matches the original cells, which carry no explicit style).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.172.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.684.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.519"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.48%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.12"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +8.26%  "
$ws.Range("E9").Value = "  +3.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0625"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.84%  "
$ws.Range("E11").Value = "  +0.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.922.90"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.701.66"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.40%  "
$ws.Range("E14").Value = "  +2.24%  "
$ws.Range("E15").Value = "  +3.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.169.86"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "235.79"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.66%  "
$ws.Range("E20").Value = "  +0.88%  "
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("E22").Value = "  +2.41%  "
$ws.Range("E23").Value = "  +3.90%  "
$ws.Range("E24").Value = "  -2.36%  "
$ws.Range("E25").Value = "  +0.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.14%  "
$ws.Range("E28").Value = "  +0.30%  "
$ws.Range("E29").Value = "  +0.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0504"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.73%  "
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("E32").Value = "  +1.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.537.52"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.37%  "
$ws.Range("E34").Value = "  +1.89%  "
$ws.Range("E35").Value = "  -2.04%  "
$ws.Range("E36").Value = "  +3.02%  "
$ws.Range("E37").Value = "  +3.11%  "
$ws.Range("E38").Value = "  -0.48%  "
$ws.Range("E39").Value = "  -0.55%  "
$ws.Range("E40").Value = "  +1.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.09"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.73"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.35%  "
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("E44").Value = "  -0.87%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.831.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.793"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.12%  "
$ws.Range("E48").Value = "  +5.66%  "
$ws.Range("E49").Value = "  +3.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.43%  "
$ws.Range("E51").Value = "  +0.01%  "
